# Applies the "Add mac address for devices in excel" edit:
#  - Sheet2: replace the old "Classroom/Test1/Test2/Test3/Sarah Office/Dorm"
#    scratch header cells with plain numbers, add two raw RSSI readings in
#    I2:I3, and fill in a third distance-estimate curve in column F
#    (rows 3:82) using the same POWER() model as columns B/C/D/E.
#  - Add a new Sheet3 at the end of the workbook with a small MAC-address /
#    device table, and make it the active sheet/tab.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Row 1/2 scratch header cleanup -----------------------------------
$ws2.Range("F1").Value = 38
$ws2.Range("H1").ClearContents()
$ws2.Range("J1").ClearContents()

$ws2.Range("F2").Value = 2.5
$ws2.Range("G2").ClearContents()
$ws2.Range("H2").ClearContents()

$ws2.Range("I2").Value = 78.32142
$ws2.Range("I3").Value = 80.50877

# --- Leftover one-off numbers/labels in column F/G ---------------------
$ws2.Range("G6").ClearContents()
$ws2.Range("G7").ClearContents()
$ws2.Range("G8").ClearContents()

# --- New distance column F: POWER(10,(ABS(A)-50.894736)/10/2.9614) ----
$ws2.Range("F3:F82").Formula = "=POWER(10,(ABS(A3)-50.894736)/10/2.9614)"

# --- New Sheet3: MAC address / device table -----------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $last)
$ws3.Name = "Sheet3"

$ws3.Range("A2").Value = "B0:91:22:F7:64:82"
$ws3.Range("C1").Value = "Dylan"
$ws3.Range("D1").Value = "Sarah"
$ws3.Range("C2").Value = 1
$ws3.Range("D2").Value = 9
$ws3.Range("E2").Value = "External"

$ws3.Range("A3").Value = "B0:91:22:F7:6a:dd"
$ws3.Range("C3").Value = 1
$ws3.Range("D3").Value = 9
$ws3.Range("E3").Value = "Internal"

$ws3.Range("A4").Value = "B0:91:22:F7:6B:1D"
$ws3.Range("C4").Value = 4
$ws3.Range("D4").Value = 6
$ws3.Range("E4").Value = "Internal"

$ws3.Range("E9").Select()
